$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# Updated bus voltage magnitudes (p.u.) for the 380 kV case (Case_0_139).
# $rows is a table: row number -> values for columns B,C,D,E,F,I,J,K,L,M,N (in order).
$rowData = @{
    2 = @(1.02, 1.026905519074363, 1.034161649376047, 1.027069699967144, 1.040269403592229, 1.02520380891933, 1.032066356648624, 1.036961917775336, 1.029890503051378, 1.04305223070588, 1.033532009254746)
    3 = @(1.02, 1.028391146146689, 1.035591938284478, 1.028341937396994, 1.042024406689247, 1.025406818414385, 1.033189186684463, 1.038199682250761, 1.030969131476395, 1.044615138637505, 1.034656433838012)
    4 = @(1.02, 1.029344449228114, 1.036507257230878, 1.029158642315391, 1.043141095417076, 1.025527838821604, 1.033907421830326, 1.038989939803724, 1.031660032118447, 1.045607134151081, 1.035375688960164)
    5 = @(1.02, 1.029743328404619, 1.036889649017363, 1.02950044509972, 1.043606066307424, 1.025576248777343, 1.034207399765616, 1.039319638050143, 1.031948819758312, 1.046019586661673, 1.03567609289853)
    6 = @(1.02, 1.029810191881047, 1.03695371382193, 1.029557745560343, 1.043683875172294, 1.025584232580588, 1.034257652573288, 1.039374848428394, 1.03199721126431, 1.046088571668833, 1.035726417070953)
    7 = @(1.02, 1.029349786462355, 1.036512376195754, 1.029163215525029, 1.043147325947278, 1.025528495360574, 1.033911437855223, 1.038994355139646, 1.031663897437656, 1.045612663314417, 1.035379710688278)
    8 = @(1.02, 1.027409265013405, 1.034647145644525, 1.027501020740754, 1.040866458685719, 1.02527456339692, 1.032447556515343, 1.037382449024848, 1.030256501106087, 1.043584446363668, 1.033913750468985)
    9 = @(1.02, 1.023927381277334, 1.031281148341646, 1.024521130930631, 1.036700321057757, 1.024747470173962, 1.029803328294307, 1.034459196853044, 1.027721614638217, 1.039860615989295, 1.031265767140578)
    10 = @(1.02, 1.021562386895395, 1.02898194250757, 1.022498880460239, 1.033820966908499, 1.024341834801052, 1.027995487332851, 1.03245280460864, 1.025993420788633, 1.03727436901979, 1.029455358837595)
    11 = @(1.02, 1.020527524682096, 1.027972800001205, 1.021614421518891, 1.032549244471438, 1.024153144712362, 1.027201625419682, 1.031569919061869, 1.02523568709829, 1.036129187870147, 1.028660369549451)
    12 = @(1.02, 1.020141469608562, 1.027595878143046, 1.021284539011665, 1.032073053776738, 1.024081079320785, 1.026905055365416, 1.031239816524591, 1.024952786681808, 1.035699947316765, 1.028363378331691)
    13 = @(1.02, 1.020224355534016, 1.027676824055171, 1.021355361769168, 1.032175372143831, 1.024096627370096, 1.026968747881103, 1.031310723036918, 1.025013535641333, 1.035792197055903, 1.028427161298056)
    14 = @(1.02, 1.020495647363068, 1.027941686236476, 1.021587181157819, 1.032509960665103, 1.024147228209773, 1.027177145611835, 1.031542676979128, 1.025212332105667, 1.036093786035864, 1.028635854977469)
    15 = @(1.02, 1.0206625779242, 1.02810459945437, 1.021729832210724, 1.032715603865367, 1.024178142480032, 1.027305320775951, 1.031685304051103, 1.025334624850448, 1.036279090224524, 1.02876421216502)
    16 = @(1.02, 1.021630836462602, 1.029048626466231, 1.022557390810035, 1.033904835616365, 1.024354081205251, 1.028047937597884, 1.032511098241596, 1.026043508203537, 1.037349831978979, 1.02950788358802)
    17 = @(1.02, 1.022235280667571, 1.029637126967469, 1.023074115849259, 1.0346440838005, 1.024460938016726, 1.028510778683315, 1.033025293033517, 1.026485630819213, 1.038014656332473, 1.02997138196088)
    18 = @(1.02, 1.02258680297693, 1.029979082916396, 1.023374664873054, 1.035072872972838, 1.024522008157025, 1.028779681783451, 1.033323856131819, 1.026742607220432, 1.038399998814426, 1.030240666934263)
    19 = @(1.02, 1.022706487760299, 1.030095460717027, 1.023477001420868, 1.035218673646211, 1.024542618682766, 1.028871191261408, 1.033425429282281, 1.026830076703427, 1.038530978986177, 1.030332306366174)
    20 = @(1.02, 1.022170537342119, 1.029574121864575, 1.023018764071806, 1.034565018381549, 1.024449603495485, 1.028461230535207, 1.032970265551757, 1.026438289195665, 1.037943579639027, 1.029921763448719)
    21 = @(1.02, 1.020415804851559, 1.027863748734215, 1.021518953791474, 1.032411538607048, 1.02413238225764, 1.027115824725502, 1.031474432273858, 1.025153831622132, 1.036005082952182, 1.028574447008444)
    22 = @(1.02, 1.019302903743108, 1.026776304312529, 1.02056810686892, 1.031035443289092, 1.023921481973563, 1.026260093611689, 1.030521426518493, 1.024337869448392, 1.034763851645161, 1.027717500658296)
    23 = @(1.02, 1.019893799756088, 1.027353937545918, 1.021072924389509, 1.031767057892356, 1.024034375725629, 1.026714675433672, 1.031027833352782, 1.024771230479314, 1.035423999922332, 1.028172728038607)
    24 = @(1.02, 1.022199795275191, 1.029602595193405, 1.023043777776201, 1.034600752059025, 1.024454728963008, 1.02848362249757, 1.032995134296453, 1.026459683636431, 1.03797570366837, 1.029944187210237)
    25 = @(1.02, 1.024835106066279, 1.032160920550728, 1.025297678039597, 1.037795092101475, 1.024893240616407, 1.030494742141802, 1.035224930517181, 1.028383583101046, 1.04084136454989, 1.031958162875043)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $r).Value = $vals[$i]
    }
}

